# Plantilla_Prediccion_Consumo.xlsx - update prediction run metadata
# - Updates the "Creado" timestamp line to a "Última predicción" line
# - Fills in the predicted consumption values for rows 6-8 with a new
#   "predicted" highlight style (italic dark-green font / light-green fill)
# - Updates the "Entrenado" timestamp in the Instrucciones sheet

$wb = $excel.ActiveWorkbook

$wsDatos = $wb.Worksheets.Item("Datos para Predicción")
$wsInstr = $wb.Worksheets.Item("Instrucciones")

# --- 1. Header timestamp: "Creado: ..." -> "Última predicción: ..." ---
$wsDatos.Range("A3").Value = "Última predicción: 2025-11-25 10:49:00"
$wsDatos.Range("A3").Font.Italic = $true
$wsDatos.Range("A3").Font.Size = 9
$wsDatos.Range("A3").Font.Color = 25600

# --- 2. Fill in predicted values for the first 3 data rows ---
$wsDatos.Range("G6").Value = 824.6900000000001
$wsDatos.Range("G7").Value = 1008.26
$wsDatos.Range("G8").Value = 824.6900000000001

$predictedRange = $wsDatos.Range("G6:G8")
$predictedRange.Interior.Color = 13561798

# --- 3. Instrucciones sheet: update "Entrenado" timestamp ---
$wsInstr.Range("A21").Value = "  - Entrenado: 2025-11-25 09:18:38"
